{"js": "// The document's sole table is a 20-row x 5-column grid of arithmetic\n// problems (e.g. \"22+7=29\"). The commit replaces the text of every cell,\n// reading the table left-to-right / top-to-bottom, with a new set of\n// problems - the row/column layout (20x5 = 100 cells) itself is unchanged,\n// only the text inside each cell differs. We therefore load the table,\n// and overwrite its `values` with the new grid, which rewrites the text\n// of the existing run in every cell while leaving all other formatting\n// (fonts, size, paragraph alignment, table/row/cell properties) untouched.\n\n// Flat list (row-major, left-to-right / top-to-bottom) of the new cell\n// text, 100 entries for the 20x5 table.\nconst newCellValues = [\n  \"34+26=60\", \"58-54=4\", \"79-74=5\", \"12+67=79\", \"81-3=78\",\n  \"32-27=5\", \"96-79=17\", \"54+15=69\", \"89-10=79\", \"3+92=95\",\n  \"1+36=37\", \"2+74=76\", \"44-36=8\", \"30+48=78\", \"85-13=72\",\n  \"28-11=17\", \"11-5=6\", \"23+40=63\", \"34+43=77\", \"14+56=70\",\n  \"85-51=34\", \"52+7=59\", \"56-7=49\", \"65-56=9\", \"3+44=47\",\n  \"10-9=1\", \"48+36=84\", \"98-45=53\", \"4+63=67\", \"66+10=76\",\n  \"54-42=12\", \"90-23=67\", \"90-72=18\", \"18+55=73\", \"41+42=83\",\n  \"10+75=85\", \"19-11=8\", \"94-58=36\", \"50+0=50\", \"90-65=25\",\n  \"17+49=66\", \"23+66=89\", \"25+15=40\", \"18+53=71\", \"99-84=15\",\n  \"80-8=72\", \"16+82=98\", \"51-7=44\", \"90-8=82\", \"15+63=78\",\n  \"75-34=41\", \"3+73=76\", \"60-1=59\", \"83-24=59\", \"32+34=66\",\n  \"95-18=77\", \"30+56=86\", \"68-13=55\", \"37+18=55\", \"86-16=70\",\n  \"33+15=48\", \"73-45=28\", \"28+40=68\", \"99-7=92\", \"93-57=36\",\n  \"93-33=60\", \"21-19=2\", \"28+0=28\", \"73+11=84\", \"40-21=19\",\n  \"11-7=4\", \"55-2=53\", \"43-0=43\", \"6+56=62\", \"73-47=26\",\n  \"29+57=86\", \"76+5=81\", \"67-34=33\", \"24-11=13\", \"84+6=90\",\n  \"35+1=36\", \"38+43=81\", \"50-42=8\", \"28+39=67\", \"13+86=99\",\n  \"77+2=79\", \"70+12=82\", \"96-66=30\", \"67-22=45\", \"27+68=95\",\n  \"58-16=42\", \"31+66=97\", \"64-32=32\", \"5+77=82\", \"4+95=99\",\n  \"60-24=36\", \"0+60=60\", \"86-59=27\", \"28+68=96\", \"66-13=53\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values[0].length;\nconst newGrid = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  newGrid.push(newCellValues.slice(r * colCount, r * colCount + colCount));\n}\n\n// Whole-table assignment rewrites each cell's text in place, keeping the\n// existing run/paragraph formatting (fonts, size, alignment) intact.\ntable.values = newGrid;\n\nawait context.sync();\n", "ps1": "# The document's sole table is a 20-row x 5-column grid of arithmetic\n# problems (e.g. \"22+7=29\"). The commit replaces the text of every cell,\n# reading the table left-to-right / top-to-bottom, with a new set of\n# problems - the row/column layout (20x5 = 100 cells) itself is unchanged,\n# only the text inside each cell differs. We therefore walk every\n# Row/Column index of Table 1 and overwrite Cell(r,c).Range.Text, which\n# rewrites the text of the existing run in place and leaves all other\n# formatting (fonts, size, paragraph alignment, table/row/cell\n# properties) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New cell text, one array per row (row-major, top-to-bottom / left-to-right).\n$newValues = @(\n  @(\"34+26=60\",\"58-54=4\",\"79-74=5\",\"12+67=79\",\"81-3=78\"),\n  @(\"32-27=5\",\"96-79=17\",\"54+15=69\",\"89-10=79\",\"3+92=95\"),\n  @(\"1+36=37\",\"2+74=76\",\"44-36=8\",\"30+48=78\",\"85-13=72\"),\n  @(\"28-11=17\",\"11-5=6\",\"23+40=63\",\"34+43=77\",\"14+56=70\"),\n  @(\"85-51=34\",\"52+7=59\",\"56-7=49\",\"65-56=9\",\"3+44=47\"),\n  @(\"10-9=1\",\"48+36=84\",\"98-45=53\",\"4+63=67\",\"66+10=76\"),\n  @(\"54-42=12\",\"90-23=67\",\"90-72=18\",\"18+55=73\",\"41+42=83\"),\n  @(\"10+75=85\",\"19-11=8\",\"94-58=36\",\"50+0=50\",\"90-65=25\"),\n  @(\"17+49=66\",\"23+66=89\",\"25+15=40\",\"18+53=71\",\"99-84=15\"),\n  @(\"80-8=72\",\"16+82=98\",\"51-7=44\",\"90-8=82\",\"15+63=78\"),\n  @(\"75-34=41\",\"3+73=76\",\"60-1=59\",\"83-24=59\",\"32+34=66\"),\n  @(\"95-18=77\",\"30+56=86\",\"68-13=55\",\"37+18=55\",\"86-16=70\"),\n  @(\"33+15=48\",\"73-45=28\",\"28+40=68\",\"99-7=92\",\"93-57=36\"),\n  @(\"93-33=60\",\"21-19=2\",\"28+0=28\",\"73+11=84\",\"40-21=19\"),\n  @(\"11-7=4\",\"55-2=53\",\"43-0=43\",\"6+56=62\",\"73-47=26\"),\n  @(\"29+57=86\",\"76+5=81\",\"67-34=33\",\"24-11=13\",\"84+6=90\"),\n  @(\"35+1=36\",\"38+43=81\",\"50-42=8\",\"28+39=67\",\"13+86=99\"),\n  @(\"77+2=79\",\"70+12=82\",\"96-66=30\",\"67-22=45\",\"27+68=95\"),\n  @(\"58-16=42\",\"31+66=97\",\"64-32=32\",\"5+77=82\",\"4+95=99\"),\n  @(\"60-24=36\",\"0+60=60\",\"86-59=27\",\"28+68=96\",\"66-13=53\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n  }\n}\n"}
